# Update scripts with new TPM values.
# 1) Remove the two rows whose Target cluster is "Resolving-Mac" (original rows 5 and 9).
# 2) Recompute the numeric columns (E:T) for the remaining rows with the new TPM-derived values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows bottom-up so row indices of the rows still to be removed do not shift.
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(5).Delete()

# New values for the surviving rows (after the deletions/shift above), columns E through T.
$newValues = @{
    2 = @(1, 0.3333333333333333, 0.112551, 0.337653, 0.7825479339666589, 0.7825479339666588, 3, 1, 0.8155003333333334, 2.446501, 0.1910612426590028, 0.1910612426590029, 0.091785378017, 0.8260684021529999, 0.1495145807039051, 0.1495145807039051)
    3 = @(1, 0.3333333333333333, 0.112551, 0.337653, 0.7825479339666589, 0.7825479339666588, 3, 1, 3.333134333333334, 9.999403000000001, 0.7809105179307759, 0.780910517930776, 0.375147602351, 3.376328421159, 0.6110999124195622, 0.6110999124195622)
    4 = @(1, 0.3333333333333333, 0.112551, 0.337653, 0.7825479339666589, 0.7825479339666588, 3, 1, 0.119632, 0.358896, 0.02802823941022116, 0.02802823941022117, 0.013464701232, 0.121182311088, 0.02193344084319146, 0.02193344084319146)
    5 = @(1, 0.3333333333333333, 0.03127533333333334, 0.09382600000000001, 0.2174520660333412, 0.2174520660333412, 3, 1, 0.8155003333333334, 2.446501, 0.1910612426590028, 0.1910612426590029, 0.02550504475844445, 0.229545402826, 0.04154666195509771, 0.04154666195509771)
    6 = @(1, 0.3333333333333333, 0.03127533333333334, 0.09382600000000001, 0.2174520660333412, 0.2174520660333412, 3, 1, 3.333134333333334, 9.999403000000001, 0.7809105179307759, 0.780910517930776, 0.1042448873197778, 0.9382039858780001, 0.1698106055112137, 0.1698106055112138)
    7 = @(1, 0.3333333333333333, 0.03127533333333334, 0.09382600000000001, 0.2174520660333412, 0.2174520660333412, 3, 1, 0.119632, 0.358896, 0.02802823941022116, 0.02802823941022117, 0.003741530677333334, 0.033673776096, 0.006094798567029707, 0.006094798567029708)
}

$cols = @(5, 6, 7, 8, 9, 10, 11, 12, 13, 14, 15, 16, 17, 18, 19, 20)  # columns E..T

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Cells.Item($row, $cols[$i]).Value2 = $vals[$i]
    }
}
